$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (header): add I4 = 2020, formatted like the other year headers (H4) ---
$ws.Range("H4").Copy()
$ws.Range("I4").PasteSpecial(-4122)
$ws.Range("I4").Value = 2020
$ws.Range("I4").VerticalAlignment = -4107

# --- Row 5: add I5 = 25.2, same style as the rest of the row (H5) ---
$ws.Range("H5").Copy()
$ws.Range("I5").PasteSpecial(-4122)
$ws.Range("I5").Value = 25.2

# --- Rows 6-13: fill in the 2020 column, same numeric style as H6:H13 but without forced vertical centering ---
$ws.Range("H6").Copy()
$ws.Range("I6").PasteSpecial(-4122)
$ws.Range("I6").Value = 39.4
$ws.Range("I6").VerticalAlignment = -4107

$ws.Range("H7").Copy()
$ws.Range("I7").PasteSpecial(-4122)
$ws.Range("I7").Value = 35.1
$ws.Range("I7").VerticalAlignment = -4107

$ws.Range("H8").Copy()
$ws.Range("I8").PasteSpecial(-4122)
$ws.Range("I8").Value = 42.7
$ws.Range("I8").VerticalAlignment = -4107

$ws.Range("H9").Copy()
$ws.Range("I9").PasteSpecial(-4122)
$ws.Range("I9").Value = 37.5
$ws.Range("I9").VerticalAlignment = -4107

$ws.Range("H10").Copy()
$ws.Range("I10").PasteSpecial(-4122)
$ws.Range("I10").Value = 40.9
$ws.Range("I10").VerticalAlignment = -4107

$ws.Range("H11").Copy()
$ws.Range("I11").PasteSpecial(-4122)
$ws.Range("I11").Value = 36.7
$ws.Range("I11").VerticalAlignment = -4107

$ws.Range("H12").Copy()
$ws.Range("I12").PasteSpecial(-4122)
$ws.Range("I12").Value = 24.7
$ws.Range("I12").VerticalAlignment = -4107

$ws.Range("H13").Copy()
$ws.Range("I13").PasteSpecial(-4122)
$ws.Range("I13").Value = -8
$ws.Range("I13").VerticalAlignment = -4107

# --- Row 14 (bottom of the table, thick bottom border): add I14 = 38.8 ---
$ws.Range("D14").Copy()
$ws.Range("I14").PasteSpecial(-4122)
$ws.Range("I14").Value = 38.8
$ws.Range("I14").NumberFormat = "0.0"
$ws.Range("I14").Font.Name = "Times New Roman"
$ws.Range("I14").Font.Size = 9
$ws.Range("I14").Font.ThemeColor = 1
$ws.Range("I14").VerticalAlignment = -4107

# --- Update the selection shown in the sheet view ---
$ws.Range("I15").Select()
